$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-28 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-29 Sunday", 2)

$d.Content.Find.Execute("867×4=3468", $true, $false, $false, $false, $false, $true, 1, $false, "987×3=2961", 2)
$d.Content.Find.Execute("817×7=5719", $true, $false, $false, $false, $false, $true, 1, $false, "312×8=2496", 2)
$d.Content.Find.Execute("688×3=2064", $true, $false, $false, $false, $false, $true, 1, $false, "689×6=4134", 2)
$d.Content.Find.Execute("937×7=6559", $true, $false, $false, $false, $false, $true, 1, $false, "801×5=4005", 2)
$d.Content.Find.Execute("755×2=1510", $true, $false, $false, $false, $false, $true, 1, $false, "832×4=3328", 2)

$d.Content.Find.Execute("286×7=2002", $true, $false, $false, $false, $false, $true, 1, $false, "970×3=2910", 2)
$d.Content.Find.Execute("160×4=640", $true, $false, $false, $false, $false, $true, 1, $false, "716×8=5728", 2)
$d.Content.Find.Execute("376×6=2256", $true, $false, $false, $false, $false, $true, 1, $false, "870×4=3480", 2)
$d.Content.Find.Execute("178×9=1602", $true, $false, $false, $false, $false, $true, 1, $false, "253×7=1771", 2)
$d.Content.Find.Execute("441×7=3087", $true, $false, $false, $false, $false, $true, 1, $false, "856×3=2568", 2)

$d.Content.Find.Execute("445×3=1335", $true, $false, $false, $false, $false, $true, 1, $false, "907×9=8163", 2)
$d.Content.Find.Execute("535×2=1070", $true, $false, $false, $false, $false, $true, 1, $false, "938×8=7504", 2)
$d.Content.Find.Execute("738×4=2952", $true, $false, $false, $false, $false, $true, 1, $false, "209×8=1672", 2)
$d.Content.Find.Execute("121×8=968", $true, $false, $false, $false, $false, $true, 1, $false, "558×6=3348", 2)
$d.Content.Find.Execute("734×2=1468", $true, $false, $false, $false, $false, $true, 1, $false, "936×6=5616", 2)

$d.Content.Find.Execute("709×9=6381", $true, $false, $false, $false, $false, $true, 1, $false, "817×6=4902", 2)
$d.Content.Find.Execute("297×3=891", $true, $false, $false, $false, $false, $true, 1, $false, "679×6=4074", 2)
$d.Content.Find.Execute("691×6=4146", $true, $false, $false, $false, $false, $true, 1, $false, "743×5=3715", 2)
$d.Content.Find.Execute("235×8=1880", $true, $false, $false, $false, $false, $true, 1, $false, "501×9=4509", 2)
$d.Content.Find.Execute("622×3=1866", $true, $false, $false, $false, $false, $true, 1, $false, "989×7=6923", 2)

$d.Content.Find.Execute("417×9=3753", $true, $false, $false, $false, $false, $true, 1, $false, "304×7=2128", 2)
$d.Content.Find.Execute("976×4=3904", $true, $false, $false, $false, $false, $true, 1, $false, "807×9=7263", 2)
$d.Content.Find.Execute("500×6=3000", $true, $false, $false, $false, $false, $true, 1, $false, "401×9=3609", 2)
$d.Content.Find.Execute("357×2=714", $true, $false, $false, $false, $false, $true, 1, $false, "496×7=3472", 2)
$d.Content.Find.Execute("449×7=3143", $true, $false, $false, $false, $false, $true, 1, $false, "164×5=820", 2)
